$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Header row (row 1)
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Data rows (rows 2-13)
$ws.Range("A2").Value = 65
$ws.Range("B2").Value = "臺灣銀行"
$ws.Range("C2").Value = "活期存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "簡東明"
$ws.Range("F2").Value = 1477218
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2012-04-24"
$ws.Range("J2").Value = "簡東明"
$ws.Range("K2").Value = 1717
$ws.Range("L2").Value = "tmpfdfe1"
$ws.Range("M2").Value = 65

$ws.Range("A3").Value = 66
$ws.Range("B3").Value = "臺灣新光商業銀行"
$ws.Range("C3").Value = "活期存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "簡東明"
$ws.Range("F3").Value = 1035741
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2012-04-24"
$ws.Range("J3").Value = "簡東明"
$ws.Range("K3").Value = 1717
$ws.Range("L3").Value = "tmpfdfe1"
$ws.Range("M3").Value = 66

$ws.Range("A4").Value = 68
$ws.Range("B4").Value = "彰化商業銀行"
$ws.Range("C4").Value = "活期存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "簡東明"
$ws.Range("F4").Value = 15913
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2012-04-24"
$ws.Range("J4").Value = "簡東明"
$ws.Range("K4").Value = 1717
$ws.Range("L4").Value = "tmpfdfe1"
$ws.Range("M4").Value = 68

$ws.Range("A5").Value = 69
$ws.Range("B5").Value = "屏東縣枋山地區農會"
$ws.Range("C5").Value = "活期存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "簡東明"
$ws.Range("F5").Value = 367026
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2012-04-24"
$ws.Range("J5").Value = "簡東明"
$ws.Range("K5").Value = 1717
$ws.Range("L5").Value = "tmpfdfe1"
$ws.Range("M5").Value = 69

$ws.Range("A6").Value = 70
$ws.Range("B6").Value = "中華郵政股份有限公司"
$ws.Range("C6").Value = "活期存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "簡東明"
$ws.Range("F6").Value = 4300
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2012-04-24"
$ws.Range("J6").Value = "簡東明"
$ws.Range("K6").Value = 1717
$ws.Range("L6").Value = "tmpfdfe1"
$ws.Range("M6").Value = 70

$ws.Range("A7").Value = 71
$ws.Range("B7").Value = "臺灣銀行屏東分行"
$ws.Range("C7").Value = "活期存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "戴錦花"
$ws.Range("F7").Value = 432537
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2012-04-24"
$ws.Range("J7").Value = "簡東明"
$ws.Range("K7").Value = 1717
$ws.Range("L7").Value = "tmpfdfe1"
$ws.Range("M7").Value = 71

$ws.Range("A8").Value = 72
$ws.Range("B8").Value = "臺灣銀行屏東分行"
$ws.Range("C8").Value = "活期存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "戴錦花"
$ws.Range("F8").Value = 113470
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2012-04-24"
$ws.Range("J8").Value = "簡東明"
$ws.Range("K8").Value = 1717
$ws.Range("L8").Value = "tmpfdfe1"
$ws.Range("M8").Value = 72

$ws.Range("A9").Value = 73
$ws.Range("B9").Value = "彰化商業銀行屏東分行"
$ws.Range("C9").Value = "活期存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "戴錦花"
$ws.Range("F9").Value = 161306
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("I9").Value = "2012-04-24"
$ws.Range("J9").Value = "簡東明"
$ws.Range("K9").Value = 1717
$ws.Range("L9").Value = "tmpfdfe1"
$ws.Range("M9").Value = 73

$ws.Range("A10").Value = 74
$ws.Range("B10").Value = "合作金庫商業銀行屏東分行"
$ws.Range("C10").Value = "活期存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "戴錦花"
$ws.Range("F10").Value = 519926
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("I10").Value = "2012-04-24"
$ws.Range("J10").Value = "簡東明"
$ws.Range("K10").Value = 1717
$ws.Range("L10").Value = "tmpfdfe1"
$ws.Range("M10").Value = 74

$ws.Range("A11").Value = 75
$ws.Range("B11").Value = "臺灣新光商業銀行東園分行"
$ws.Range("C11").Value = "活期存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "戴錦花"
$ws.Range("F11").Value = 242
$ws.Range("G11").Value = "deposit"
$ws.Range("H11").Value = "normal"
$ws.Range("I11").Value = "2012-04-24"
$ws.Range("J11").Value = "簡東明"
$ws.Range("K11").Value = 1717
$ws.Range("L11").Value = "tmpfdfe1"
$ws.Range("M11").Value = 75

$ws.Range("A12").Value = 76
$ws.Range("B12").Value = "屏東六塊厝郵局(第21支局）"
$ws.Range("C12").Value = "活期存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "戴錦花"
$ws.Range("F12").Value = 287235
$ws.Range("G12").Value = "deposit"
$ws.Range("H12").Value = "normal"
$ws.Range("I12").Value = "2012-04-24"
$ws.Range("J12").Value = "簡東明"
$ws.Range("K12").Value = 1717
$ws.Range("L12").Value = "tmpfdfe1"
$ws.Range("M12").Value = 76

$ws.Range("A13").Value = 77
$ws.Range("B13").Value = "屏東縣枋山地區農會信用部"
$ws.Range("C13").Value = "活期存款"
$ws.Range("D13").Value = "新臺幣"
$ws.Range("E13").Value = "戴錦花"
$ws.Range("F13").Value = 29295
$ws.Range("G13").Value = "deposit"
$ws.Range("H13").Value = "normal"
$ws.Range("I13").Value = "2012-04-24"
$ws.Range("J13").Value = "簡東明"
$ws.Range("K13").Value = 1717
$ws.Range("L13").Value = "tmpfdfe1"
$ws.Range("M13").Value = 77
